$d = $word.ActiveDocument

$old = "Conoscere l’importanza del “Regolamento sulla privacy” (Privacy Policy) che i servizi digitali predispongono per informare gli utenti sull’utilizzo dei dati personali raccolti con focus particolare sui social network e la profilazione degli utenti. Conoscere e applicare le misure di sicurezza, protezione, tutela della riservatezza. Proteggere i dispositivi e i contenuti e comprendere i rischi e le minacce presenti negli ambienti digitali;"
$new = "Conoscere i principali documenti italiani ed europei per la regolamentazione dell’intelligenza artificiale, le motivazioni che hanno portato a tali documentazioni, legate alla storia dell’intelligenza artificiale, al suo funzionamento, ai suoi problemi (anche ambientali) e limiti;"

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
